# Scenarios.xlsx edit script
# Adds a new "PopulationId" column (inserted before "ModelParameterSheets")
# and a new row describing a population-based scenario ("PopulationScenario").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# Insert a new column before the current column C ("ModelParameterSheets")
# so the new "PopulationId" column lands there, shifting the remaining
# columns one to the right.
$ws.Columns.Item(3).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 3).Value = "PopulationId"
$ws.Cells.Item(1, 3).Style = $ws.Cells.Item(1, 2).Style

# Add a new row describing a population scenario.
$ws.Cells.Item(4, 1).Value = "PopulationScenario"
$ws.Cells.Item(4, 2).Value = "Indiv"
$ws.Cells.Item(4, 3).Value = "TestPopulation"
$ws.Cells.Item(4, 4).Value = "Global"
$ws.Cells.Item(4, 5).Value = "Aciclovir_iv_250mg"
$ws.Cells.Item(4, 6).Value = 12
$ws.Cells.Item(4, 7).Value = "h"
$ws.Cells.Item(4, 8).Value = $false
$ws.Cells.Item(4, 11).Value = "Aciclovir.pkml"

# Resize the new column to fit its content, like the other data columns.
$ws.Columns.Item(3).AutoFit()

$ws.Range("I4:J4").Select()
